$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 2; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "*All*") {
        $cell.Value2 = "<b>All</b>"
    }
}
